$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 5 new rows (8-12) by copying row 7 formatting, then filling values below
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(8).Insert() | Out-Null
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(9).Insert() | Out-Null
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(10).Insert() | Out-Null
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(11).Insert() | Out-Null
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(12).Insert() | Out-Null
$excel.CutCopyMode = 0

# Step 2: set cell values + row heights for rows 2-12

# Row 2
$ws.Range("A2").Value = 'CP_GESCLSERDOM_001'
$ws.Range("B2").Value = 'ingreso a la vista Autodiagnostico'
$ws.Range("C2").Value = 'Positivo'
$ws.Range("D2").Value = 'eCenter'
$ws.Range("E2").Value = 'El usuario debe tener permisos para acceder a la vista'
$ws.Range("F2").Value = '1. Clic en módulo eCenter
2. Scroll en el contenedor de aplicaciones
3. Clic en "Gestion clientes y servicios domiciliarios"'
$ws.Range("G2").Value = 'ID de cliente válido'
$ws.Range("H2").Value = 'El sistema debe redirigido correctamente la vista Gestion clientes y servicios domiciliarios'
$ws.Range("I2").Value = 'La vista Autodiagnóstico se cargó sin errores'
$ws.Range("J2").Value = 'OK'
$ws.Range("K2").Value = 'SI'
$ws.Range("L2").Value = 'N/A'
$ws.Rows.Item(2).RowHeight = 98.25

# Row 3
$ws.Range("A3").Value = 'CP_GESCLSERDOM_002'
$ws.Range("B3").Value = 'Filtro de búsqueda cliente por ID_DEAL'
$ws.Range("C3").Value = 'Positivo'
$ws.Range("D3").Value = 'eCenter'
$ws.Range("E3").Value = 'El ID_DEAL del cliente debe existir en base de datos'
$ws.Range("F3").Value = '1.Abrir modal de filtros
2.Desplegar select de filtros.
3.Diligenciar el campo de ID_DEAL
4.Clic en "Aplicar filtros"'
$ws.Range("G3").Value = 'ID de cliente válido'
$ws.Range("H3").Value = 'El sistema debe mostrar la  informacion del cliente'
$ws.Range("I3").Value = 'se visualizó la información del cliente correctamente'
$ws.Range("J3").Value = 'OK'
$ws.Range("K3").Value = 'SI'
$ws.Range("L3").Value = 'N/A'
$ws.Rows.Item(3).RowHeight = 70.5

# Row 4
$ws.Range("A4").Value = 'CP_GESCLSERDOM_003'
$ws.Range("B4").Value = 'Ver información técnica asociada'
$ws.Range("C4").Value = 'Positivo'
$ws.Range("D4").Value = 'eCenter'
$ws.Range("E4").Value = 'Cliente en estado ACTIVO y tipo RESIDENCIAL debe existir en la base de datos'
$ws.Range("F4").Value = '1. Seleccionar cliente con estado ACTIVO y tipo RESIDENCIAL.
2. Clic en Opciones.
3. Clic en Ver información técnica asociada'
$ws.Range("G4").Value = 'ID de cliente válido'
$ws.Range("H4").Value = 'Se muestra un modal con la información técnica del cliente.'
$ws.Range("I4").Value = 'El modal de información técnica se abrió y mostró los datos correctos del cliente.'
$ws.Range("J4").Value = 'OK'
$ws.Range("K4").Value = 'SI'
$ws.Range("L4").Value = 'N/A'
$ws.Rows.Item(4).RowHeight = 90.75

# Row 5
$ws.Range("A5").Value = 'CP_GESCLSERDOM_004'
$ws.Range("B5").Value = 'Reconfiguración 
del cliente'
$ws.Range("C5").Value = 'Positivo'
$ws.Range("D5").Value = 'eCenter'
$ws.Range("E5").Value = 'Cliente seleccionado y visible en la tabla'
$ws.Range("F5").Value = '1. Clic en Opciones.
2. Seleccionar Reconfiguración.
3. Clic en Reconfigurar.
4. Confirmar en el modal con Sí.'
$ws.Range("G5").Value = 'ID de cliente válido'
$ws.Range("H5").Value = 'Se inicia el proceso de reconfiguración y se muestran barras de progreso.'
$ws.Range("I5").Value = 'El proceso de reconfiguración inició y las barras de progreso se mostraron correctamente.'
$ws.Range("J5").Value = 'OK'
$ws.Range("K5").Value = 'SI'
$ws.Range("L5").Value = 'N/A'
$ws.Rows.Item(5).RowHeight = 66.75

# Row 6
$ws.Range("A6").Value = 'CP_GESCLSERDOM_005'
$ws.Range("B6").Value = 'Ver dispositivos del cliente'
$ws.Range("C6").Value = 'Positivo'
$ws.Range("D6").Value = 'eCenter'
$ws.Range("E6").Value = 'Cliente en estado ACTIVO seleccionado'
$ws.Range("F6").Value = '1. Clic en Opciones.
2. Seleccionar opción "Suspensión"'
$ws.Range("G6").Value = 'ID de cliente válido'
$ws.Range("H6").Value = 'El sistema abre un modal que lista los dispositivos asociados al cliente.'
$ws.Range("I6").Value = 'El modal se abrió y mostró la lista de dispositivos asociados.'
$ws.Range("J6").Value = 'OK'
$ws.Range("K6").Value = 'SI'
$ws.Range("L6").Value = 'N/A'
$ws.Rows.Item(6).RowHeight = 56.25

# Row 7
$ws.Range("A7").Value = 'CP_GESCLSERDOM_006'
$ws.Range("B7").Value = 'Ver y enviar documentos (Acta de instalación y Contrato)'
$ws.Range("C7").Value = 'Positivo'
$ws.Range("D7").Value = 'eCenter'
$ws.Range("E7").Value = 'Cliente con 
documentos asociados'
$ws.Range("F7").Value = '1. Clic en Opciones.
2. Seleccionar Ver documentos.
3. Para cada documento:
 a. Clic Ver documento.
 b. Clic Enviar al correo.
 c. Clic Descargar.'
$ws.Range("G7").Value = 'N/A'
$ws.Range("H7").Value = 'El sistema abre un modal para visualizar los documentos del cliente(Acta de instalación y Contrato)'
$ws.Range("I7").Value = 'Acta y Contrato se visualizaron, enviaron y descargaron según lo esperado.'
$ws.Range("J7").Value = 'OK'
$ws.Range("K7").Value = 'SI'
$ws.Range("L7").Value = 'N/A'
$ws.Rows.Item(7).RowHeight = 94.5

# Row 8
$ws.Range("A8").Value = 'CP_GESCLSERDOM_007'
$ws.Range("B8").Value = 'Ver detalle del proceso'
$ws.Range("C8").Value = 'Positivo'
$ws.Range("D8").Value = 'eCenter'
$ws.Range("E8").Value = 'Cliente seleccionado'
$ws.Range("F8").Value = '1. Clic en Opciones.
2. Seleccionar Detalle del proceso.'
$ws.Range("G8").Value = 'ID de cliente válido'
$ws.Range("H8").Value = 'Se despliega un modal con el historial y detalle de procesos del cliente.'
$ws.Range("I8").Value = 'El modal se abrió y mostró el historial y detalle de procesos correctamente.'
$ws.Range("J8").Value = 'OK'
$ws.Range("K8").Value = 'SI'
$ws.Range("L8").Value = 'N/A'
$ws.Rows.Item(8).RowHeight = 88.5

# Row 9
$ws.Range("A9").Value = 'CP_GESCLSERDOM_008'
$ws.Range("B9").Value = 'Suspensión del cliente'
$ws.Range("C9").Value = 'Positivo'
$ws.Range("D9").Value = 'eCenter'
$ws.Range("E9").Value = 'Cliente en estado ACTIVO y con plan vigente.'
$ws.Range("F9").Value = '1. Seleccionar cliente.
2. Clic en Opciones.
3. Seleccionar Suspensión de servicio.
4. Confirmar en el modal con botón Sí.'
$ws.Range("G9").Value = 'ID de cliente válido'
$ws.Range("H9").Value = 'El sistema inicia la suspensión, muestra barra de progreso y cambia el estado del cliente a SUSPENDIDO.'
$ws.Range("I9").Value = 'Suspensión ejecutada y estado del cliente actualizado a SUSPENDIDO.'
$ws.Range("J9").Value = 'OK'
$ws.Range("K9").Value = 'SI'
$ws.Range("L9").Value = 'N/A'
$ws.Rows.Item(9).RowHeight = 96

# Row 10
$ws.Range("A10").Value = 'CP_GESCLSERDOM_009'
$ws.Range("B10").Value = 'Renexion de servicio suspendido'
$ws.Range("C10").Value = 'Positivo'
$ws.Range("D10").Value = 'eCenter'
$ws.Range("E10").Value = 'Cliente previamente suspendido (Estado = SUSPENDIDO).'
$ws.Range("F10").Value = '1. Seleccionar cliente.
2. Clic en Opciones.
3. Seleccionar Reconexión de servicio.
4. Confirmar en el modal con botón Sí.'
$ws.Range("G10").Value = 'ID de cliente suspendido'
$ws.Range("H10").Value = 'El sistema reconecta el servicio, muestra barra de progreso y cambia el estado del cliente a ACTIVO.'
$ws.Range("I10").Value = 'Servicio re-conectado y estado del cliente actualizado a ACTIVO.'
$ws.Range("J10").Value = 'OK'
$ws.Range("K10").Value = 'SI'
$ws.Range("L10").Value = 'N/A'
$ws.Rows.Item(10).RowHeight = 99

# Row 11
$ws.Range("A11").Value = 'CP_GESCLSERDOM_009'
$ws.Range("B11").Value = 'Renexion de servicio suspendido'
$ws.Range("C11").Value = 'Positivo'
$ws.Range("D11").Value = 'eCenter'
$ws.Range("E11").Value = 'Cliente previamente suspendido (Estado = SUSPENDIDO).'
$ws.Range("F11").Value = '1. Seleccionar cliente.
2. Clic en Opciones.
3. Seleccionar Reconexión de servicio.
4. Confirmar en el modal con botón Sí.'
$ws.Range("G11").Value = 'ID de cliente suspendido'
$ws.Range("H11").Value = 'El sistema reconecta el servicio, muestra barra de progreso y cambia el estado del cliente a ACTIVO.'
$ws.Range("I11").Value = 'Servicio re-conectado y estado del cliente actualizado a ACTIVO.'
$ws.Range("J11").Value = 'OK'
$ws.Range("K11").Value = 'SI'
$ws.Range("L11").Value = 'N/A'
$ws.Rows.Item(11).RowHeight = 100.5

# Row 12
$ws.Range("A12").Value = 'CP_GESCLSERDOM_010'
$ws.Range("B12").Value = 'Cambio de plan de servicio'
$ws.Range("C12").Value = 'Positivo'
$ws.Range("D12").Value = 'eCenter'
$ws.Range("E12").Value = 'Cliente en estado ACTIVO y con al menos un plan disponible para cambio.'
$ws.Range("F12").Value = '1. Seleccionar cliente.
2. Clic en Opciones.
3. Seleccionar Cambio de plan.
4. En el modal, escoger nuevo plan de la lista.
5. Clic en Guardar/Confirmar cambio.'
$ws.Range("G12").Value = 'ID de cliente y nombre del nuevo plan'
$ws.Range("H12").Value = 'El sistema actualiza el plan contratado y muestra confirmación de cambio exitoso.'
$ws.Range("I12").Value = 'Cambio de plan ejecutado correctamente y plan actualizado.'
$ws.Range("J12").Value = 'OK'
$ws.Range("K12").Value = 'SI'
$ws.Range("L12").Value = 'N/A'
$ws.Rows.Item(12).RowHeight = 120.75

# Step 3: selection
$ws.Range("I3").Select() | Out-Null

Write-Host "done"